# "Fruta / hortaliza, semanal" — weekly refresh of the Brócoli /
# Vega Monumental Concepción subset: a new week's record is inserted
# at the top of the data block (row 151), pushing the existing rows
# (151-172) down by one (to 152-173).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 151; Excel shifts rows 151:172 down to 152:173
# and carries the row-above formatting (incl. the date style on D)
# onto the new row automatically.
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row with this week's reading.
$ws.Cells.Item(151, 1).Value  = 11
$ws.Cells.Item(151, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(151, 3).Value  = "Bíobío"
$ws.Cells.Item(151, 4).Value  = 44505
$ws.Cells.Item(151, 5).Value  = 8
$ws.Cells.Item(151, 6).Value  = 100112023
$ws.Cells.Item(151, 7).Value  = "Brócoli"
$ws.Cells.Item(151, 8).Value  = "Sin especificar"
$ws.Cells.Item(151, 9).Value  = "Primera"
$ws.Cells.Item(151, 10).Value = 2700
$ws.Cells.Item(151, 11).Value = 550
$ws.Cells.Item(151, 12).Value = 600
$ws.Cells.Item(151, 13).Value = 578
$ws.Cells.Item(151, 14).Value = "$/unidad"
$ws.Cells.Item(151, 15).Value = "Región del Maule"
$ws.Cells.Item(151, 16).Value = 578
$ws.Cells.Item(151, 17).Value = 1
$ws.Cells.Item(151, 18).Value = "Hortaliza"
